$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 6-46) holds running totals that were re-baselined: every
# existing value drops by 30 (62..102 -> 32..72) while column A/B stay put.
for ($r = 6; $r -le 46; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 - 30
}

# A few wrapped-text rows had an explicit (stale) row height; re-autofit so
# Excel drops back to the sheet's default row height for them.
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).AutoFit()
$ws.Rows.Item(23).AutoFit()

# Scroll the window down toward the bottom of the table and move the
# selection to the new working cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 36
$win.ScrollColumn = 1
$ws.Range("G43").Select()
